# "Updated script to accomodate for multiple series"
#
# The authored change removes the empty, unused "Content Placeholder 4"
# placeholder shape (shape id 5) from the third slide (sldId 828) of the
# deck. (PowerPoint's own save machinery also re-stamps the cached
# datetimeFigureOut date fields in the handout/notes masters and rewrites
# collaboration-tracking parts such as revisionInfo.xml / changesInfo1.xml
# whenever the file is opened and re-saved in the authoring app - those are
# incidental side effects of the save pipeline, not an explicit edit, so
# they are not reproduced here.)

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(3)

for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.Name -eq "Content Placeholder 4") {
        $shape.Delete()
    }
}
